$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Previously these cells held numeric placeholder values (999 / 100).
# Item values are now read line by line as strings, so correct them
# to proper string values consistent with the rest of the column.
$ws.Range("A1").Value = "item1"
$ws.Range("B10").Value = "value 10"

# Move the active selection to B11, past the last data row.
$ws.Range("B11").Select()
